$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" both contain the same event table and both
# need the same new row inserted (the diff shows an identical change
# applied to both worksheets).
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Helper: write a plain-text value into a cell without Excel's
    # automatic date/number recognition kicking in (used for the
    # "开始时间" column which holds literal yyyy-mm-dd strings).
    function Set-TextValue($range, $value) {
        $range.NumberFormat = "@"
        $range.Value = $value
        $range.Style = "Normal"
    }

    # Column A cells (A1:A.. ) all share the same bold/centered/bordered
    # style. Use A1 as the canonical source so copying its format never
    # introduces a brand new style entry.
    $ws.Range("A1").Copy() | Out-Null

    # Current row 6 (AEO event) needs to move down to row 7, keeping its
    # data but bumping the running counter in column A from 5 to 6.
    Set-TextValue $ws.Range("B7") "2024-08-17"
    $ws.Range("A7").Value = 6
    $ws.Range("C7").Value = "丽水·AEO纯白礼赞动漫嘉年华"
    $ws.Range("D7").Value = "城北街1001号 爱依·时尚婚宴中心"
    $ws.Range("E7").Value = "2024.08.17 09:00-08.17 16:00"
    $ws.Range("F7").Value = 161
    $ws.Range("G7").Value = 55
    $ws.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=86779"
    $ws.Range("I7").Value = "//i2.hdslb.com/bfs/openplatform/202406/MxJ3oNjt1717405405850.jpeg"
    $ws.Range("A7").PasteSpecial(-4122)

    # Current row 5 (CCAC 七夕 event) needs to move down to row 6, data
    # unchanged, counter in column A stays 5.
    Set-TextValue $ws.Range("B6") "2024-08-10"
    $ws.Range("A6").Value = 5
    $ws.Range("C6").Value = "丽水·CCAC动漫七夕（回馈展）"
    $ws.Range("D6").Value = "中东路848号(解放街交汇) 飞达国际大酒店"
    $ws.Range("E6").Value = "2024.08.10 09:00-08.10 17:00"
    $ws.Range("F6").Value = 33
    $ws.Range("G6").Value = 29.9
    $ws.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=86567"
    $ws.Range("I6").Value = "//i0.hdslb.com/bfs/openplatform/202405/tsOzbBRx1717015539538.png"
    $ws.Range("A6").PasteSpecial(-4122)

    # Row 5 becomes the newly added event (丽水·樱卡动漫游戏嘉年华),
    # counter in column A stays 4, style on A5 is already correct.
    Set-TextValue $ws.Range("B5") "2024-08-03"
    $ws.Range("A5").Value = 4
    $ws.Range("C5").Value = "丽水·樱卡动漫游戏嘉年华"
    $ws.Range("D5").Value = "中东路848号(解放街交汇) 飞达国际大酒店"
    $ws.Range("E5").Value = "2024.08.03 10:00-08.03 17:00"
    $ws.Range("F5").Value = 0
    $ws.Range("G5").Value = 50
    $ws.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=87276"
    $ws.Range("I5").Value = "//i0.hdslb.com/bfs/openplatform/202406/bVp0Zg1B1718172430380.jpeg"
}

Write-Host "edit complete"
